# Applies revisions described in the commit:
#   "Implementation revisions to budget by model number and modify model
#    output written to Excel (i.e. oh, rates, asset info, formulas vs hard code)"
#
# Concretely, on the "Vacation Summary" sheet:
#   - C3 becomes a formula that averages C2 and C4 instead of a hard-coded value
#   - F3 (the "prod_days" hard input) is revised from 37 to 32
#   - the active selection on "Vacation Summary" moves to A6
#   - the active selection on "Vac Entitlement Table" moves to E19

$wb = $excel.ActiveWorkbook

$wsVacSummary = $wb.Worksheets.Item("Vacation Summary")
$wsVacEntitlement = $wb.Worksheets.Item("Vac Entitlement Table")

# Replace the hard-coded C3 value with a formula averaging C2 and C4
$wsVacSummary.Range("C3").Formula = "=(C2+C4)/2"

# Revise the hard-coded prod_days input in F3
$wsVacSummary.Range("F3").Value = 32

# Update the saved selections on each sheet to match where the user left off.
# Select the "Vac Entitlement Table" cell first, then finish on the
# "Vacation Summary" cell so that sheet stays the active/selected tab,
# matching the original workbook state.
$wsVacEntitlement.Range("E19").Select()
$wsVacSummary.Range("A6").Select()

$wb.Save()
